# Update "想去人数" (want-to-go headcount) figures scraped for the
# gh-pages data refresh (commit "Update gh-pages to output generated at 456a3b4").
#
# Sheet "展览" and sheet "全部类型" both list the same events (just at
# different row offsets), so each sheet gets its own set of F-column
# (column 6) updates per the refreshed scrape.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value  = 784    # was 783
$wsExhibit.Range("F3").Value  = 60     # was 59
$wsExhibit.Range("F7").Value  = 159    # was 158
$wsExhibit.Range("F8").Value  = 347    # was 346
$wsExhibit.Range("F9").Value  = 457    # was 456
$wsExhibit.Range("F10").Value = 515    # was 516
$wsExhibit.Range("F12").Value = 11812  # was 11790
$wsExhibit.Range("F13").Value = 5419   # was 5418

# --- Sheet "全部类型" ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value  = 784    # was 783
$wsAll.Range("F3").Value  = 60     # was 59
$wsAll.Range("F9").Value  = 159    # was 158
$wsAll.Range("F10").Value = 347    # was 346
$wsAll.Range("F11").Value = 457    # was 456
$wsAll.Range("F12").Value = 515    # was 516
$wsAll.Range("F14").Value = 11813  # was 11790
$wsAll.Range("F16").Value = 5419   # was 5418
